$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value2 = 1352375.2
$ws.Range("I62").Value2 = 2061802.4
$ws.Range("J62").Value2 = 169996.67
$ws.Range("K62").Value2 = 2061802.4
$ws.Range("L62").Value2 = 169996.67
$ws.Range("M62").Value2 = -2061178.4
$ws.Range("N62").Value2 = -171244.67
$ws.Range("H65").Value2 = 1352375.2
$ws.Range("I65").Value2 = 2061802.4
$ws.Range("J65").Value2 = 169996.67
$ws.Range("K65").Value2 = 10309012
$ws.Range("L65").Value2 = 849983.3500000001
$ws.Range("M65").Value2 = -10305892
$ws.Range("N65").Value2 = -856223.3500000001
$ws.Range("H106").Value2 = 1000000
$ws.Range("I106").Value2 = 0
$ws.Range("K106").Value2 = 0
$ws.Range("M106").ClearContents()
$ws.Range("H112").Value2 = 64055.5
$ws.Range("J112").Value2 = 1900.75
$ws.Range("L112").Value2 = 5702.25
$ws.Range("N112").Value2 = -7918.25
$ws.Range("H138").Value2 = 2803.0625
$ws.Range("J138").Value2 = 2970.375
$ws.Range("L138").Value2 = 8911.125
$ws.Range("N138").Value2 = -19191.125

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 1060.1666
$ws.Range("I2").Value2 = 1098.8889
$ws.Range("J2").Value2 = 944
$ws.Range("K2").Value2 = 1098.8889
$ws.Range("L2").Value2 = 944
$ws.Range("M2").Value2 = -985.8888999999999
$ws.Range("N2").Value2 = -1170
$ws.Range("H5").Value2 = 6775
$ws.Range("I5").Value2 = 50
$ws.Range("K5").Value2 = 50
$ws.Range("M5").Value2 = 62
$ws.Range("H88").Value2 = 15154146
$ws.Range("I88").Value2 = 33334322
$ws.Range("J88").Value2 = 3998.8333
$ws.Range("K88").Value2 = 33334322
$ws.Range("L88").Value2 = 3998.8333
$ws.Range("M88").Value2 = -33333916
$ws.Range("N88").Value2 = -4810.8333
$ws.Range("H91").Value2 = 15154146
$ws.Range("I91").Value2 = 33334322
$ws.Range("J91").Value2 = 3998.8333
$ws.Range("K91").Value2 = 33334322
$ws.Range("L91").Value2 = 3998.8333
$ws.Range("M91").Value2 = -33332918
$ws.Range("N91").Value2 = -6806.8333
$ws.Range("H116").Value2 = 1060.1666
$ws.Range("I116").Value2 = 1098.8889
$ws.Range("J116").Value2 = 944
$ws.Range("K116").Value2 = 1098.8889
$ws.Range("L116").Value2 = 944
$ws.Range("M116").Value2 = 1195.1111
$ws.Range("N116").Value2 = -5532
$ws.Range("H131").Value2 = 62999.5
$ws.Range("J131").Value2 = 62999.5
$ws.Range("L131").Value2 = 62999.5
$ws.Range("N131").Value2 = -73079.5
$ws.Range("H132").Value2 = 40001970
$ws.Range("I132").Value2 = 45456380
$ws.Range("K132").Value2 = 136369140
$ws.Range("M132").Value2 = -136366610

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 1060.1666
$ws.Range("I3").Value2 = 1098.8889
$ws.Range("J3").Value2 = 944
$ws.Range("K3").Value2 = 1098.8889
$ws.Range("L3").Value2 = 944
$ws.Range("M3").Value2 = -984.8888999999999
$ws.Range("N3").Value2 = -1172
$ws.Range("H4").Value2 = 6775
$ws.Range("I4").Value2 = 50
$ws.Range("K4").Value2 = 50
$ws.Range("M4").Value2 = 65
$ws.Range("H35").Value2 = 69091.58
$ws.Range("J35").Value2 = 69091.58
$ws.Range("L35").Value2 = 69091.58
$ws.Range("N35").Value2 = -69711.58
$ws.Range("H86").Value2 = 4622.7144
$ws.Range("I86").Value2 = 5042.8
$ws.Range("K86").Value2 = 5042.8
$ws.Range("M86").Value2 = -3919.8
$ws.Range("H89").Value2 = 4622.7144
$ws.Range("I89").Value2 = 5042.8
$ws.Range("K89").Value2 = 25214
$ws.Range("M89").Value2 = -19598
$ws.Range("H94").Value2 = 2655.2222
$ws.Range("I94").Value2 = 1539.4
$ws.Range("K94").Value2 = 1539.4
$ws.Range("M94").Value2 = -1088.4
$ws.Range("H99").Value2 = 1849.75
$ws.Range("I99").Value2 = 1849.75
$ws.Range("K99").Value2 = 1849.75
$ws.Range("M99").Value2 = -351.75
$ws.Range("H125").Value2 = 0
$ws.Range("J125").Value2 = 0
$ws.Range("L125").ClearContents()
$ws.Range("N125").Value2 = 0
$ws.Range("H134").Value2 = 2589.2896
$ws.Range("I134").Value2 = 2500
$ws.Range("K134").Value2 = 7500
$ws.Range("M134").Value2 = -4965

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value2 = 55779.5
$ws.Range("J20").Value2 = 55779.5
$ws.Range("L20").Value2 = 55779.5
$ws.Range("N20").Value2 = -56251.5
$ws.Range("H26").Value2 = 962
$ws.Range("I26").Value2 = 962
$ws.Range("J26").Value2 = 0
$ws.Range("K26").Value2 = 962
$ws.Range("L26").Value2 = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value2 = -675
$ws.Range("H30").Value2 = 55779.5
$ws.Range("J30").Value2 = 55779.5
$ws.Range("L30").Value2 = 55779.5
$ws.Range("N30").Value2 = -55961.5
$ws.Range("H62").Value2 = 5236.8184
$ws.Range("I62").Value2 = 3367.5
$ws.Range("K62").Value2 = 3367.5
$ws.Range("M62").Value2 = -2743.5
$ws.Range("H65").Value2 = 5236.8184
$ws.Range("I65").Value2 = 3367.5
$ws.Range("K65").Value2 = 16837.5
$ws.Range("M65").Value2 = -13717.5
$ws.Range("H128").Value2 = 55779.5
$ws.Range("J128").Value2 = 55779.5
$ws.Range("L128").Value2 = 55779.5
$ws.Range("N128").Value2 = -65739.5
$ws.Range("H130").Value2 = 50000
$ws.Range("J130").Value2 = 50000
$ws.Range("L130").Value2 = 50000
$ws.Range("N130").Value2 = -60040
$ws.Range("H132").Value2 = 2631.8125
$ws.Range("I132").Value2 = 2489.5925
$ws.Range("K132").Value2 = 7468.7775
$ws.Range("M132").Value2 = -4938.7775

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value2 = 645.8889
$ws.Range("J107").Value2 = 496.14285
$ws.Range("L107").Value2 = 1488.42855
$ws.Range("N107").Value2 = -5328.428550000001
$ws.Range("H132").Value2 = 1540.8334
$ws.Range("J132").Value2 = 2500
$ws.Range("L132").Value2 = 22500
$ws.Range("N132").Value2 = -27560
$ws.Range("H134").Value2 = 1872.3636
$ws.Range("I134").Value2 = 1872.3636
$ws.Range("K134").Value2 = 5617.0908
$ws.Range("M134").Value2 = -547.0907999999999
$ws.Range("H137").Value2 = 1365.8334
$ws.Range("I137").Value2 = 998.75
$ws.Range("K137").Value2 = 2996.25
$ws.Range("M137").Value2 = 2103.75

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value2 = 2688.1667
$ws.Range("I80").Value2 = 2478
$ws.Range("K80").Value2 = 2478
$ws.Range("M80").Value2 = -1480
$ws.Range("H83").Value2 = 2688.1667
$ws.Range("I83").Value2 = 2478
$ws.Range("K83").Value2 = 12390
$ws.Range("M83").Value2 = -7398
$ws.Range("H130").Value2 = 66799.8
$ws.Range("J130").Value2 = 66799.8
$ws.Range("L130").Value2 = 66799.8
$ws.Range("N130").Value2 = -76839.8
$ws.Range("H132").Value2 = 3808.7188
$ws.Range("I132").Value2 = 3089.8096
$ws.Range("K132").Value2 = 9269.4288
$ws.Range("M132").Value2 = -6739.4288

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value2 = 12404.167
$ws.Range("I68").Value2 = 7625
$ws.Range("J68").Value2 = 14793.75
$ws.Range("K68").Value2 = 7625
$ws.Range("L68").Value2 = 14793.75
$ws.Range("M68").Value2 = -6876
$ws.Range("N68").Value2 = -16291.75
$ws.Range("H71").Value2 = 12404.167
$ws.Range("I71").Value2 = 7625
$ws.Range("J71").Value2 = 14793.75
$ws.Range("K71").Value2 = 38125
$ws.Range("L71").Value2 = 73968.75
$ws.Range("M71").Value2 = -34381
$ws.Range("N71").Value2 = -81456.75
$ws.Range("H132").Value2 = 3337.8918
$ws.Range("I132").Value2 = 3204.3215
$ws.Range("J132").Value2 = 3753.4443
$ws.Range("K132").Value2 = 9612.9645
$ws.Range("L132").Value2 = 11260.3329
$ws.Range("M132").Value2 = -7082.9645
$ws.Range("N132").Value2 = -16320.3329

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value2 = 3406287.8
$ws.Range("I62").Value2 = 3972502.5
$ws.Range("J62").Value2 = 9000
$ws.Range("K62").Value2 = 3972502.5
$ws.Range("L62").Value2 = 9000
$ws.Range("M62").Value2 = -3971878.5
$ws.Range("N62").Value2 = -10248
$ws.Range("H65").Value2 = 3406287.8
$ws.Range("I65").Value2 = 3972502.5
$ws.Range("J65").Value2 = 9000
$ws.Range("K65").Value2 = 19862512.5
$ws.Range("L65").Value2 = 45000
$ws.Range("M65").Value2 = -19859392.5
$ws.Range("N65").Value2 = -51240
$ws.Range("H81").Value2 = 8339083.5
$ws.Range("I81").Value2 = 3236.1333
$ws.Range("K81").Value2 = 6472.2666
$ws.Range("M81").Value2 = -5411.2666
$ws.Range("H84").Value2 = 8339083.5
$ws.Range("I84").Value2 = 3236.1333
$ws.Range("K84").Value2 = 32361.333
$ws.Range("M84").Value2 = -27057.333
$ws.Range("H122").Value2 = 3855.7856
$ws.Range("I122").Value2 = 3873.5
$ws.Range("K122").Value2 = 11620.5
$ws.Range("M122").Value2 = -9170.5
$ws.Range("H136").Value2 = 3644.1538
$ws.Range("I136").Value2 = 2431.6924
$ws.Range("K136").Value2 = 7295.0772
$ws.Range("M136").Value2 = -4745.0772
